# Applies the BSR item re-shuffle / Grand Total restatement described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 8 ----
$ws.Range("A8").Value = "Each"
$ws.Range("C8").Value = 19
$ws.Range("D8").Value = "'3.0"
$ws.Range("E8").Value = 'P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F8").Value = 23
$ws.Range("G8").Value = "'437.00"

# ---- Row 9 ----
$ws.Range("A9").Value = ""
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = "'11.0"
$ws.Range("E9").Value = "S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = "'0.00"

# ---- Row 10 ----
$ws.Range("A10").Value = "Mtr."
$ws.Range("C10").Value = 74
$ws.Range("D10").Value = "'19"
$ws.Range("E10").Value = "2 x 2.5 sq. mm. + 1x1.5sqmm"
$ws.Range("F10").Value = 81
$ws.Range("G10").Value = "'5994.00"

# ---- Row 11 ----
$ws.Range("A11").Value = "Set"
$ws.Range("C11").Value = 98
$ws.Range("D11").Value = "'13.0"
$ws.Range("E11").Value = 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F11").Value = 5733
$ws.Range("G11").Value = "'561834.00"

# ---- Row 12 ----
$ws.Range("A12").Value = ""
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = "'38"
$ws.Range("E12").Value = "Grand Total"
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = "'0.00"

# ---- Row 14 (Grand Total Rs.) ----
$ws.Range("G14").Value = "'568265.00"
$ws.Range("H14").Value = "'568265.00"

# ---- Row 16 (NET PAYABLE AMOUNT Rs.) ----
$ws.Range("G16").Value = "'568265.00"
$ws.Range("H16").Value = "'568265.00"
